# Generate Report for Handback
# Refresh the handback timestamps for the file that was just handed back
# (f889c0fd-9ec9-4624-b0c8-bc6f7f90dd34.md), across the Overview sheet and
# the per-locale (zh-cn / de-de) detail sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
# Row 3 corresponds to f889c0fd-9ec9-4624-b0c8-bc6f7f90dd34.md
# Column G = "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-23 18:51:46"

# --- zh-cn sheet ------------------------------------------------------
# Row 3 corresponds to f889c0fd-9ec9-4624-b0c8-bc6f7f90dd34.md
# Column H = "Correspond Handoff Datetime"
# Column K = "Correspond Handback DateTime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-23 18:51:41"
$wsZhCn.Range("K3").Value = "2016-08-23 18:51:59"

# --- de-de sheet ------------------------------------------------------
# Row 3 corresponds to f889c0fd-9ec9-4624-b0c8-bc6f7f90dd34.md
# Column H = "Correspond Handoff Datetime"
# Column K = "Correspond Handback DateTime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-23 18:51:46"
$wsDeDe.Range("K3").Value = "2016-08-23 18:52:18"
